# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect the latest generated data (commit 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 31
$ws1.Range("F4").Value = 1395
$ws1.Range("F5").Value = 320
$ws1.Range("F7").Value = 10729
$ws1.Range("F9").Value = 83
$ws1.Range("F11").Value = 1039
$ws1.Range("F12").Value = 715
$ws1.Range("F13").Value = 12060
$ws1.Range("F14").Value = 12513

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 31
$ws4.Range("F5").Value = 1395
$ws4.Range("F6").Value = 320
$ws4.Range("F8").Value = 10729
$ws4.Range("F10").Value = 83
$ws4.Range("F12").Value = 1039
$ws4.Range("F13").Value = 715
$ws4.Range("F14").Value = 12060
$ws4.Range("F15").Value = 12513
